# Add a new portfolio item "Full Petential" as row 7 of the items sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("items")

# Widen column F (showcased) to fit the new, longer title text.
# (75.45 is the closest COM ColumnWidth input that rounds to the target
# stored OOXML width of 76.33203125 chars.)
$ws.Columns.Item(6).ColumnWidth = 75.45

# New row 7 data. Cells are written in the same order the shared-string
# table grows (title, application, development, asset_folder, showcased,
# media) so new <si> entries land at the expected indices.
$ws.Range("A7").Value = "padded"
$ws.Range("B7").Value = 6
$ws.Range("C7").Value = 6
$ws.Range("D7").Value = "FULL PETENTIAL"
$ws.Range("E7").Value = "Identity"
$ws.Range("H7").Value = "Identity, Packaging"
$ws.Range("I7").Value = "<p>Objective: An internship turned into a paid gig! Created company identity. Assets included illustrative logo, type logo, color platform, mailers, and packaging of canine and feline treats.</p>`n<p>Packaged samples available.</p>"
$ws.Range("K7").Value = "6.Full-Petential"
$ws.Range("F7").Value = "Full Petential Identity"
$ws.Range("G7").Value = "Print, Packaging"

# The multi-paragraph text in I7 otherwise triggers an auto row-height
# bump; every data row in this sheet keeps the fixed 18pt custom height.
$ws.Rows.Item(7).RowHeight = 18
